$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title column: strip the surrounding quote marks that were baked into the
# original strings. Row 2 keeps a trailing apostrophe and is entered with a
# leading apostrophe so Excel records it as quote-prefixed text.
$ws.Range("D2").Value = "'Branch Manager'"
$ws.Range("D3").Value = "Assistant Branch Manager"
$ws.Range("D4").Value = "Loan Officer"
$ws.Range("D5").Value = "Teller"
$ws.Range("D6").Value = "Teller"
$ws.Range("D7").Value = "Teller"
$ws.Range("D8").Value = "Teller"
$ws.Range("D9").Value = "Teller"
$ws.Range("D10").Value = "Personal Banker"
$ws.Range("D11").Value = "Personal Banker"
$ws.Range("D12").Value = "Personal Banker"
$ws.Range("D13").Value = "Customer Service Representative"
$ws.Range("D14").Value = "Financial Advisor"
$ws.Range("D15").Value = "Financial Advisor"
$ws.Range("D16").Value = "Financial Advisor"

# Is_Admin column: was stored as boolean TRUE/FALSE, now plain 1/0 numbers.
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0

# Move the active selection to D2, matching the saved cursor position.
$ws.Range("D2").Select() | Out-Null
